$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tri = $s.Shapes.Item(2)
$tri.Width = 2396067 / 12700.0
Write-Output $tri.Width
